$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a plain number-like string (e.g. "313.00").
# Excel would normally auto-convert these to numeric values, losing the
# original text formatting (trailing zeros, exact digit count). Force the
# Price column cells to Text format first so the literal string is preserved,
# exactly like the source data feed (which stores prices as text strings).
$textPriceCells = @(
    "D5"
    "D7"
    "D8"
    "D9"
    "D10"
    "D11"
    "D13"
    "D14"
    "D15"
    "D16"
    "D17"
    "D18"
    "D21"
    "D23"
    "D24"
    "D25"
    "D26"
    "D27"
    "D29"
    "D30"
    "D31"
    "D32"
    "D33"
    "D34"
    "D35"
    "D36"
    "D37"
    "D39"
    "D40"
    "D41"
    "D42"
    "D43"
    "D44"
    "D46"
    "D47"
    "D48"
    "D49"
    "D50"
    "D51"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (column D) and Volume(1h) (column E) values.
$ws.Range("D2").Value = '27.825.63'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '1.906.75'
$ws.Range("D5").Value = '313.00'
$ws.Range("E5").Value = '  +0.26%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.5223'
$ws.Range("E7").Value = '  +7.27%  '
$ws.Range("D8").Value = '0.3788'
$ws.Range("D9").Value = '0.07234'
$ws.Range("E9").Value = '  -1.24%  '
$ws.Range("D10").Value = '0.9122'
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").Value = '21.28'
$ws.Range("E11").Value = '  +3.53%  '
$ws.Range("D12").Value = '1.946.05'
$ws.Range("E12").Value = '  +2.71%  '
$ws.Range("D13").Value = '0.07648'
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("D14").Value = '5.454'
$ws.Range("E14").Value = '  -0.53%  '
$ws.Range("D15").Value = '92.23'
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").Value = '0.000008702'
$ws.Range("E17").Value = '  -0.86%  '
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '27.852.82'
$ws.Range("E19").Value = '  +0.69%  '
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").Value = '5.156'
$ws.Range("E21").Value = '  +0.71%  '
$ws.Range("D22").Value = '2.161.85'
$ws.Range("E22").Value = '  +2.27%  '
$ws.Range("D23").Value = '10.87'
$ws.Range("E23").Value = '  +1.13%  '
$ws.Range("D24").Value = '6.639'
$ws.Range("E24").Value = '  +0.45%  '
$ws.Range("D25").Value = '153.48'
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("D26").Value = '1.870'
$ws.Range("E26").Value = '  -2.09%  '
$ws.Range("D27").Value = '2.171'
$ws.Range("E27").Value = '  +0.76%  '
$ws.Range("D29").Value = '114.86'
$ws.Range("E29").Value = '  -0.58%  '
$ws.Range("D30").Value = '4.866'
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("D31").Value = '0.09019'
$ws.Range("E31").Value = '  +1.39%  '
$ws.Range("D32").Value = '4.868'
$ws.Range("E32").Value = '  +5.13%  '
$ws.Range("D33").Value = '3.178'
$ws.Range("E33").Value = '  -0.74%  '
$ws.Range("D34").Value = '1.239'
$ws.Range("E34").Value = '  +1.42%  '
$ws.Range("D35").Value = '0.7804'
$ws.Range("E35").Value = '  +1.87%  '
$ws.Range("D36").Value = '0.02095'
$ws.Range("E36").Value = '  +2.73%  '
$ws.Range("D37").Value = '2.613'
$ws.Range("E38").Value = '  +3.31%  '
$ws.Range("D39").Value = '0.5575'
$ws.Range("E39").Value = '  +1.89%  '
$ws.Range("D40").Value = '1.092'
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("D41").Value = '0.05285'
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("D42").Value = '6.729'
$ws.Range("E42").Value = '  -2.32%  '
$ws.Range("D43").Value = '116.17'
$ws.Range("E43").Value = '  +3.78%  '
$ws.Range("D44").Value = '8.559'
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("D46").Value = '0.4814'
$ws.Range("E46").Value = '  +0.62%  '
$ws.Range("D47").Value = '10.52'
$ws.Range("E47").Value = '  -1.24%  '
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("D49").Value = '1.622'
$ws.Range("E49").Value = '  -0.79%  '
$ws.Range("D50").Value = '66.98'
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("D51").Value = '0.05992'
$ws.Range("E51").Value = '  -0.96%  '
